$d = $word.ActiveDocument

# --- Locate the paragraph that needs the typo fix -------------------------
# (the one that, before editing, contains "with two except")
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*with two except*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the paragraph containing 'with two except'."
}

# --- 1) Fix the typo -------------------------------------------------------
# "...source data table with two except that t..." -> "...source data table except that t..."
$findRange = $target.Range
$replaced = $findRange.Find.Execute("with two except", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "except", 2)
if (-not $replaced) {
    throw "Could not find/replace 'with two except'."
}

# --- 2) Italicize just the word "names" -----------------------------------
# "...the column names are not changed..." -> "...the column *names* are not changed..."
$locateRange = $target.Range
$located = $locateRange.Find.Execute("column names are not changed", $true, $false, $false, $false, `
                                      $false, $true, 1, $false, "", 0)
if (-not $located) {
    throw "Could not find 'column names are not changed' after fixing the typo."
}

$namesStart = $locateRange.Start + "column ".Length
$namesEnd = $namesStart + "names".Length
$rNames = $d.Range($namesStart, $namesEnd)
if ($rNames.Text -ne "names") {
    throw "Offset computation for the word 'names' landed on the wrong text: '$($rNames.Text)'."
}
$rNames.Font.Italic = $true
